# Insert a new column "id_i" / 42 at the very left of the sheet (before
# the current column A), shifting the existing columns
# (fecha_im/active_energy_im/active_power_im and their data) one
# position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before column A; this shifts the existing cell
# content (values + styles) from A,B,C to B,C,D respectively.
$ws.Columns.Item(1).Insert()

# Populate the new first column with the header and the new id value.
$ws.Range("A1").Value = "id_i"
$ws.Range("A2").Value = 42

# Match the saved selection state (active cell A2).
$ws.Range("A2").Select()

Write-Output "done"
